$d = $word.ActiveDocument

# Step 1: merge the three runs "- " / "Calculate square" / " root." into a
# single run of text "- Calculate square root." by doing a plain text
# replace (this collapses the run split because Find/Replace writes the
# replacement into one run).
$d.Content.Find.Execute("- Calculate square root.", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "- Calculate square root.", 2)

# Step 2: highlight that exact phrase in yellow.
$rng = $d.Content
$rng.Find.Execute("- Calculate square root.", $true, $true, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$rng.HighlightColorIndex = 7
